$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 inherits A1's existing (fill+font) style before we touch values, so both
# header cells land on the same style record instead of minting two new ones.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Bold both header cells together so the engine reuses a single new font/style
# for the pair instead of minting one per cell.
$ws.Range("A1:B1").Font.Bold = $true

# Header text: "ID" in A1, "Unit Name" in B1. Write B1 first so the shared
# string table picks up "Unit Name" (index 0) ahead of "ID" (index 1).
$ws.Range("B1").Value = "Unit Name"
$ws.Range("A1").Value = "ID"

# Column widths (39.140625 / 47 characters, expressed as the nearest COM
# ColumnWidth inputs this engine's pixel-snapping rounds back to those values)
$ws.Columns("A").ColumnWidth = 38.35
$ws.Columns("B").ColumnWidth = 46.15

# Restore the selection to A6, as captured in the saved view state.
$ws.Range("A6").Select() | Out-Null
